# Update the cputest latency table with newly measured benchmark numbers
# (more PDFs included in the run; a bug around cloudflare n=40/45 rows is
# still being tracked separately).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.8272290939215802
$ws.Range("E3").Value = 0.6919655038133441

$ws.Range("C4").Value = 3.373424
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = 55

$ws.Range("D5").Value = 24
$ws.Range("E5").Value = 3.379368066236442

$ws.Range("D6").Value = 483
$ws.Range("E6").Value = 31

$ws.Range("D7").Value = 251

$ws.Range("C8").Value = 357
$ws.Range("D8").Value = 590
$ws.Range("E8").Value = 75

$ws.Range("C9").Value = 1155
$ws.Range("D9").Value = 1195
$ws.Range("E9").Value = 75

$ws.Range("C10").Value = 3867
$ws.Range("D10").Value = 4355
$ws.Range("E10").Value = 211

$ws.Range("C11").Value = 2340
$ws.Range("D11").Value = 2569

$ws.Range("C13").Value = 13062
$ws.Range("D13").Value = 13396
$ws.Range("E13").Value = 799

$ws.Range("C14").Value = 43038
$ws.Range("D14").Value = 44632
$ws.Range("E14").Value = 2044

$ws.Range("C15").Value = 25812
$ws.Range("D15").Value = 28017
$ws.Range("E15").Value = 675
